$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 42602.009479166663
$ws.Range("B20").Value = "Named"
$ws.Range("C20").Value = 2741
$ws.Range("D20").Value = 120
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 60
$ws.Range("M20").Value = 40

$ws.Range("A21").Value = 42602.014409722222
$ws.Range("B21").Value = "Named"
$ws.Range("C21").Value = 3602
$ws.Range("D21").Value = 120
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 60
$ws.Range("M21").Value = 40

$ws.Range("A22").Value = 42602.01525462963
$ws.Range("B22").Value = "Named"
$ws.Range("C22").Value = 2894
$ws.Range("D22").Value = 120
$ws.Range("E22").Value = 5
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 60
$ws.Range("M22").Value = 40

$ws.Range("A23").Value = 42602.495185185187
$ws.Range("B23").Value = "Named"
$ws.Range("C23").Value = 3179
$ws.Range("D23").Value = 151
$ws.Range("E23").Value = 8
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 66
$ws.Range("I23").Value = 33
$ws.Range("J23").Value = 3
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 49
$ws.Range("M23").Value = 49

$ws.Range("A24").Value = 42602.495462962965
$ws.Range("B24").Value = "Named"
$ws.Range("C24").Value = 2956
$ws.Range("D24").Value = 151
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 66
$ws.Range("I24").Value = 33
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 49
$ws.Range("M24").Value = 49

$ws.Range("A25").Value = 42602.495717592596
$ws.Range("B25").Value = "Named"
$ws.Range("C25").Value = 2574
$ws.Range("D25").Value = 151
$ws.Range("E25").Value = 8
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 66
$ws.Range("I25").Value = 33
$ws.Range("J25").Value = 3
$ws.Range("K25").Value = 3
$ws.Range("L25").Value = 49
$ws.Range("M25").Value = 49

$ws.Range("A26").Value = 42602.495787037034
$ws.Range("B26").Value = "Named"
$ws.Range("C26").Value = 1886
$ws.Range("D26").Value = 151
$ws.Range("E26").Value = 8
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 66
$ws.Range("I26").Value = 33
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 3
$ws.Range("L26").Value = 49
$ws.Range("M26").Value = 49
